$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Repayment schedule")

# Insert a new (blank) column before column N to make room for an
# additional data column in the repayment schedule.
$ws.Columns("N:N").Insert()
$ws.Columns("N:N").ColumnWidth = 10.7109375

# Make the "Repayment schedule" sheet the active tab and move the
# selection to L14, matching the saved workbook view state.
$ws.Activate()
$ws.Range("L14").Select()
